# Applies the "Draft project plan" update:
#  - workbook: set calcMode to manual
#  - Gantt chart sheet: reset view (no frozen/topLeft offset, no cached selection),
#    re-style the "Containgency" label, split the big merged helper cell so
#    the leftover rows get their own (different) formatting
#  - Details sheet: rename a task, add a new "Specification" detail row with
#    two extra columns, insert a blank row before the "Implementation" block,
#    and set the print/page setup for that sheet
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Workbook level: switch calculation to manual
# ---------------------------------------------------------------------------
$wb.Application.Calculation = -4135   # xlCalculationManual

# ---------------------------------------------------------------------------
# 2. Sheet1 "Gantt chart"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Gantt chart")

# Reset the view: no frozen top-left offset, no stored selection
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null

# Re-centre the "Containgency" label (keeps its existing green fill)
$contingency = $ws1.Range("A11")
$contingency.HorizontalAlignment = -4108   # xlCenter

# The helper merge used to span J5:K12; it now only spans J5:K11, so unmerge
# first and re-merge the smaller block, giving the freed-up J12:K12 cells
# their own (unfilled, wrap + vertically centred) formatting
$ws1.Range("J5:K12").UnMerge() | Out-Null
$ws1.Range("J5:K11").Merge() | Out-Null

$leftover = $ws1.Range("J12:K12")
$leftover.Interior.Pattern = -4142   # xlNone
$leftover.HorizontalAlignment = -4142  # xlNone / general
$leftover.VerticalAlignment = -4108    # xlCenter
$leftover.WrapText = $true

# ---------------------------------------------------------------------------
# 3. Sheet2 "Details for each"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Details for each")

# Rename the UI design task
$ws2.Range("B11").Value = "App layout design (UI) Storyboard"

# Insert a new blank row before the old row 16 ("Implementation" section)
$ws2.Rows.Item(16).Insert() | Out-Null

# Shift the "Specification" detail down into row 12 and give it two new
# supporting columns
$ws2.Range("B12").Value = "Specification"
$ws2.Range("C12").Value = "requirements"
$ws2.Range("D12").Value = "criteria"

# Page setup for this sheet
$ws2.PageSetup.PaperSize = 9     # xlPaperA4
$ws2.PageSetup.Orientation = 1   # xlPortrait

$ws2.Activate() | Out-Null
$ws2.Range("B13").Select() | Out-Null
